# chore: update Sheets via scheduled runner
# Refresh cached market-board price/profit figures (columns H-N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 742.8333
$ws.Range("I2").Value = 733.2
$ws.Range("K2").Value = 733.2
$ws.Range("M2").Value = -620.2

$ws.Range("H15").Value = 651.9091
$ws.Range("I15").Value = 651.9091
$ws.Range("K15").Value = 1955.7273
$ws.Range("M15").Value = -1786.7273

$ws.Range("H17").Value = 474.94446
$ws.Range("J17").Value = 474.94446
$ws.Range("L17").Value = 1424.83338
$ws.Range("N17").Value = -1760.83338

$ws.Range("H29").Value = 59
$ws.Range("I29").Value = 59
$ws.Range("K29").Value = 177
$ws.Range("M29").Value = 104

$ws.Range("H38").Value = 368.64285
$ws.Range("I38").Value = 409.1
$ws.Range("J38").Value = 267.5
$ws.Range("K38").Value = 1227.3
$ws.Range("L38").Value = 802.5
$ws.Range("M38").Value = -855.3000000000002
$ws.Range("N38").Value = -1546.5

$ws.Range("H70").Value = 4625
$ws.Range("J70").Value = 5884.2
$ws.Range("L70").Value = 17652.6
$ws.Range("N70").Value = -18192.6

$ws.Range("H73").Value = 4625
$ws.Range("J73").Value = 5884.2
$ws.Range("L73").Value = 17652.6
$ws.Range("N73").Value = -19524.6

$ws.Range("H132").Value = 279958.72
$ws.Range("I132").Value = 2131.2258
$ws.Range("K132").Value = 6393.6774
$ws.Range("M132").Value = -3863.6774

$ws.Range("H135").Value = 742.2432
$ws.Range("I135").Value = 746.25
$ws.Range("J135").Value = 716.6
$ws.Range("K135").Value = 6716.25
$ws.Range("L135").Value = 6449.400000000001
$ws.Range("M135").Value = -4181.25
$ws.Range("N135").Value = -11519.4

$ws.Range("H138").Value = 2320.4666
$ws.Range("I138").Value = 2604.2666
$ws.Range("K138").Value = 7812.7998
$ws.Range("M138").Value = -2672.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1561.5
$ws.Range("I2").Value = 1640.4231
$ws.Range("K2").Value = 1640.4231
$ws.Range("M2").Value = -1527.4231

$ws.Range("H32").Value = 6059.8286
$ws.Range("I32").Value = 2579.2144
$ws.Range("K32").Value = 2579.2144
$ws.Range("M32").Value = -2292.2144

$ws.Range("H61").Value = 3460.3157
$ws.Range("I61").Value = 2696.5
$ws.Range("J61").Value = 5599
$ws.Range("K61").Value = 2696.5
$ws.Range("L61").Value = 5599
$ws.Range("M61").Value = -2484.5
$ws.Range("N61").Value = -6023

$ws.Range("H63").Value = 5832
$ws.Range("I63").Value = 2231.6667
$ws.Range("K63").Value = 2231.6667
$ws.Range("M63").Value = -1545.6667

$ws.Range("H66").Value = 5832
$ws.Range("I66").Value = 2231.6667
$ws.Range("K66").Value = 11158.3335
$ws.Range("M66").Value = -7726.333500000001

$ws.Range("H116").Value = 1561.5
$ws.Range("I116").Value = 1640.4231
$ws.Range("K116").Value = 1640.4231
$ws.Range("M116").Value = 653.5769

$ws.Range("H122").Value = 2174
$ws.Range("I122").Value = 2199.5
$ws.Range("K122").Value = 6598.5
$ws.Range("M122").Value = -4148.5

$ws.Range("H132").Value = 2699.3572
$ws.Range("I132").Value = 1739.4445
$ws.Range("K132").Value = 5218.333500000001
$ws.Range("M132").Value = -2688.333500000001

$ws.Range("H136").Value = 3460.3157
$ws.Range("I136").Value = 2696.5
$ws.Range("J136").Value = 5599
$ws.Range("K136").Value = 8089.5
$ws.Range("L136").Value = 16797
$ws.Range("M136").Value = -5539.5
$ws.Range("N136").Value = -21897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1561.5
$ws.Range("I3").Value = 1640.4231
$ws.Range("K3").Value = 1640.4231
$ws.Range("M3").Value = -1526.4231

$ws.Range("H107").Value = 2749.5
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80

$ws.Range("H134").Value = 3216.8572
$ws.Range("I134").Value = 2686.3635
$ws.Range("K134").Value = 8059.0905
$ws.Range("M134").Value = -5524.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 285
$ws.Range("I22").Value = 213.5
$ws.Range("J22").Value = 499.5
$ws.Range("K22").Value = 213.5
$ws.Range("L22").Value = 499.5
$ws.Range("M22").Value = 136.5
$ws.Range("N22").Value = -1199.5

$ws.Range("H31").Value = 3034.3
$ws.Range("I31").Value = 2545.4443
$ws.Range("J31").Value = 3434.2727
$ws.Range("K31").Value = 2545.4443
$ws.Range("L31").Value = 3434.2727
$ws.Range("M31").Value = -2250.4443
$ws.Range("N31").Value = -4024.2727

$ws.Range("H34").Value = 3034.3
$ws.Range("I34").Value = 2545.4443
$ws.Range("J34").Value = 3434.2727
$ws.Range("K34").Value = 2545.4443
$ws.Range("L34").Value = 3434.2727
$ws.Range("M34").Value = -2343.4443
$ws.Range("N34").Value = -3838.2727

$ws.Range("H94").Value = 2550.6365
$ws.Range("I94").Value = 3633.5
$ws.Range("J94").Value = 1251.2
$ws.Range("K94").Value = 3633.5
$ws.Range("L94").Value = 1251.2
$ws.Range("M94").Value = -3182.5
$ws.Range("N94").Value = -2153.2

$ws.Range("H122").Value = 2460.7646
$ws.Range("I122").Value = 2455.3
$ws.Range("J122").Value = 2468.5715
$ws.Range("K122").Value = 7365.900000000001
$ws.Range("L122").Value = 7405.7145
$ws.Range("M122").Value = -4915.900000000001
$ws.Range("N122").Value = -12305.7145

$ws.Range("H132").Value = 3284.7942
$ws.Range("I132").Value = 3261.36
$ws.Range("K132").Value = 9784.08
$ws.Range("M132").Value = -7254.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1166
$ws.Range("I22").Value = 799.2
$ws.Range("K22").Value = 2397.6
$ws.Range("M22").Value = -2228.6

$ws.Range("H24").Value = 236.75
$ws.Range("I24").Value = 236.75
$ws.Range("K24").Value = 710.25
$ws.Range("M24").Value = -480.25

$ws.Range("H27").Value = 1166
$ws.Range("I27").Value = 799.2
$ws.Range("K27").Value = 2397.6
$ws.Range("M27").Value = -2295.6

$ws.Range("H34").Value = 370.3
$ws.Range("I34").Value = 175
$ws.Range("J34").Value = 663.25
$ws.Range("K34").Value = 525
$ws.Range("L34").Value = 1989.75
$ws.Range("M34").Value = -441
$ws.Range("N34").Value = -2157.75

$ws.Range("H39").Value = 5601.875
$ws.Range("J39").Value = 6944.6665
$ws.Range("L39").Value = 20833.9995
$ws.Range("N39").Value = -21421.9995

$ws.Range("H55").Value = 3255.1428
$ws.Range("I55").Value = 1400
$ws.Range("J55").Value = 3997.2
$ws.Range("K55").Value = 4200
$ws.Range("L55").Value = 11991.6
$ws.Range("M55").Value = -4023
$ws.Range("N55").Value = -12345.6

$ws.Range("H117").Value = 833
$ws.Range("J117").Value = 1100
$ws.Range("L117").Value = 3300
$ws.Range("N117").Value = -10184

$ws.Range("H128").Value = 176690
$ws.Range("I128").Value = 176690
$ws.Range("K128").Value = 530070
$ws.Range("M128").Value = -525090

$ws.Range("H131").Value = 37471.13
$ws.Range("J131").Value = 1857.125
$ws.Range("L131").Value = 5571.375
$ws.Range("N131").Value = -15651.375

$ws.Range("H132").Value = 2308.8948
$ws.Range("I132").Value = 2133.5
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 19201.5
$ws.Range("L132").Value = 25200
$ws.Range("M132").Value = -16671.5
$ws.Range("N132").Value = -30260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 49997
$ws.Range("J109").Value = 49997
$ws.Range("L109").Value = 49997
$ws.Range("N109").Value = -52077

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2409.2856
$ws.Range("I46").Value = 1777.3846
$ws.Range("J46").Value = 3436.125
$ws.Range("K46").Value = 1777.3846
$ws.Range("L46").Value = 3436.125
$ws.Range("M46").Value = -1589.3846
$ws.Range("N46").Value = -3812.125

$ws.Range("H50").Value = 61000
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H132").Value = 3337.353
$ws.Range("I132").Value = 3232.5557
$ws.Range("K132").Value = 9697.667099999999
$ws.Range("M132").Value = -7167.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 17498.5
$ws.Range("I37").Value = 14998
$ws.Range("K37").Value = 14998
$ws.Range("M37").Value = -14795

$ws.Range("H132").Value = 5630.6665
$ws.Range("J132").Value = 2865
$ws.Range("L132").Value = 8595
$ws.Range("N132").Value = -13655
